# Apply updated crypto symbol list values (Mon Feb 13 15:43:26 UTC 2023 refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold numeric/percent-looking text; force text format
# so Excel keeps them as strings instead of converting to numbers.
$textCells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "E23", "D24", "E24", "D25", "E25", "D26", "D38", "E38", "D39", "E39", "D40", "E40", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "E48", "D49", "E49", "D50", "E50", "D51", "E51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "290.77"
$ws.Range("E2").Value = "-6.17%"
$ws.Range("D3").Value = "40.06"
$ws.Range("E3").Value = "-2.71%"
$ws.Range("D4").Value = "5.042"
$ws.Range("E4").Value = "-3.04%"
$ws.Range("D5").Value = "0.07327"
$ws.Range("E5").Value = "-4.71%"
$ws.Range("D7").Value = "1.549"
$ws.Range("E7").Value = "-9.03%"
$ws.Range("D8").Value = "0.9116"
$ws.Range("E8").Value = "-3.04%"
$ws.Range("D9").Value = "0.1196"
$ws.Range("E9").Value = "-5.77%"
$ws.Range("D10").Value = "0.1743"
$ws.Range("E10").Value = "-5.01%"
$ws.Range("D11").Value = "0.08682"
$ws.Range("E11").Value = "-4.79%"
$ws.Range("D12").Value = "0.04158"
$ws.Range("E12").Value = "-1.98%"
$ws.Range("D13").Value = "0.1052"
$ws.Range("E13").Value = "0.04%"
$ws.Range("D14").Value = "0.001274"
$ws.Range("E14").Value = "-1.26%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.005816"
$ws.Range("E15").Value = "-1.29%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.395"
$ws.Range("E16").Value = "1.33%"
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").Value = "2.397"
$ws.Range("E17").Value = "-1.16%"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").Value = "0.3285"
$ws.Range("E18").Value = "-1.05%"
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D19").Value = "7.570"
$ws.Range("E19").Value = "1.35%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "0.1352"
$ws.Range("E20").Value = "-0.04%"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "0.2883"
$ws.Range("E21").Value = "5.92%"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "0.03848"
$ws.Range("E22").Value = "-4.05%"
$ws.Range("E23").Value = "0.13%"
$ws.Range("D24").Value = "0.003892"
$ws.Range("E24").Value = "-8.14%"
$ws.Range("D25").Value = "0.0001280"
$ws.Range("E25").Value = "0.68%"
$ws.Range("D26").Value = "0.0003728"
$ws.Range("D38").Value = "0.02334"
$ws.Range("E38").Value = "-7.92%"
$ws.Range("D39").Value = "0.05020"
$ws.Range("E39").Value = "-5.74%"
$ws.Range("D40").Value = "0.007660"
$ws.Range("E40").Value = "-2.23%"
$ws.Range("E41").Value = "163.04%"
$ws.Range("D42").Value = "0.1273"
$ws.Range("E42").Value = "-3.15%"
$ws.Range("D43").Value = "0.007373"
$ws.Range("E43").Value = "10.54%"
$ws.Range("D44").Value = "0.006967"
$ws.Range("E44").Value = "-6.16%"
$ws.Range("D45").Value = "0.3155"
$ws.Range("E45").Value = "2.17%"
$ws.Range("D46").Value = "0.00006515"
$ws.Range("E46").Value = "-4.00%"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.09%"
$ws.Range("E48").Value = "14.03%"
$ws.Range("D49").Value = "0.004206"
$ws.Range("E49").Value = "35.53%"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").Value = "-0.09%"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").Value = "-0.09%"
